$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.005673666666666667
$ws.Range("H2").Value = 0.017021
$ws.Range("I2").Value = 0.003624806445217209
$ws.Range("J2").Value = 0.003624806445217209
$ws.Range("M2").Value = 0.00535
$ws.Range("N2").Value = 0.01605
$ws.Range("O2").Value = 0.003591913026022235
$ws.Range("P2").Value = 0.003591913026022235
$ws.Range("Q2").Value = 0.00003035411666666667
$ws.Range("R2").Value = 0.00027318705
$ws.Range("S2").Value = 0.00001301998948738505
$ws.Range("T2").Value = 0.00001301998948738505
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.005673666666666667
$ws.Range("H3").Value = 0.017021
$ws.Range("I3").Value = 0.003624806445217209
$ws.Range("J3").Value = 0.003624806445217209
$ws.Range("O3").Value = 0.9964080869739778
$ws.Range("P3").Value = 0.9964080869739778
$ws.Range("Q3").Value = 0.008420328415666667
$ws.Range("R3").Value = 0.075782955741
$ws.Range("S3").Value = 0.003611786455729824
$ws.Range("T3").Value = 0.003611786455729824
$ws.Range("I4").Value = 0.8811233428718493
$ws.Range("J4").Value = 0.8811233428718493
$ws.Range("M4").Value = 0.00535
$ws.Range("N4").Value = 0.01605
$ws.Range("O4").Value = 0.003591913026022235
$ws.Range("P4").Value = 0.003591913026022235
$ws.Range("Q4").Value = 0.007378523833333333
$ws.Range("R4").Value = 0.06640671449999999
$ws.Range("S4").Value = 0.003164918412793652
$ws.Range("T4").Value = 0.003164918412793652
$ws.Range("I5").Value = 0.8811233428718493
$ws.Range("J5").Value = 0.8811233428718493
$ws.Range("O5").Value = 0.9964080869739778
$ws.Range("P5").Value = 0.9964080869739778
$ws.Range("S5").Value = 0.8779584244590557
$ws.Range("T5").Value = 0.8779584244590557
$ws.Range("I6").Value = 0.1152518506829335
$ws.Range("J6").Value = 0.1152518506829335
$ws.Range("M6").Value = 0.00535
$ws.Range("N6").Value = 0.01605
$ws.Range("O6").Value = 0.003591913026022235
$ws.Range("P6").Value = 0.003591913026022235
$ws.Range("Q6").Value = 0.0009651185999999999
$ws.Range("R6").Value = 0.008686067399999999
$ws.Range("S6").Value = 0.0004139746237411984
$ws.Range("T6").Value = 0.0004139746237411984
$ws.Range("I7").Value = 0.1152518506829335
$ws.Range("J7").Value = 0.1152518506829335
$ws.Range("O7").Value = 0.9964080869739778
$ws.Range("P7").Value = 0.9964080869739778
$ws.Range("S7").Value = 0.1148378760591923
$ws.Range("T7").Value = 0.1148378760591923
